# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns for the 8a3a5b66... / b4c9410e... entries (rows 4 & 5) on both the
# zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-23 08:24:31"
$wsZh.Range("E5").Value = "2016-03-23 08:24:31"
$wsZh.Range("H4").Value = "2016-03-23 08:24:54"
$wsZh.Range("H5").Value = "2016-03-23 08:24:54"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-23 08:24:35"
$wsDe.Range("E5").Value = "2016-03-23 08:24:35"
$wsDe.Range("H4").Value = "2016-03-23 08:25:01"
$wsDe.Range("H5").Value = "2016-03-23 08:25:01"
